# Planificación.xlsx — "Add files via upload" edit
# Updates the project-planner activity list: several activity rows are
# renamed/re-dated to reflect a new set of experiment tasks, and the
# highlighted-period selector and view are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title (unchanged text, but touched by the original edit) ---
$ws.Range("B1").Value = "Proyecto de título"

# --- Highlighted period selector ---
$ws.Range("H2").Value = 18

# --- Activity table (rows 5-19): ACTIVIDAD / INICIO / DURACION / INICIO REAL / DURACION REAL / % COMPLETADO ---

# Row 5
$ws.Range("B5").Value = "Obtención DB"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1

# Row 6
$ws.Range("B6").Value = "Planificación de tiempos"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1

# Row 7
$ws.Range("B7").Value = "Revision bibliografica"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 0.9

# Row 8
$ws.Range("B8").Value = "Analizar DB"
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 1

# Row 9
$ws.Range("B9").Value = "Elección de videoclips"
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1

# Row 10
$ws.Range("B10").Value = "Selección de frames"
$ws.Range("C10").Value = 6
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 1

# Row 11
$ws.Range("B11").Value = "Preparación de datos"
$ws.Range("C11").Value = 7
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 1

# Row 12
$ws.Range("B12").Value = "Experimentración con modelos"
$ws.Range("C12").Value = 9
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 8
$ws.Range("F12").Value = 5
$ws.Range("G12").Value = 0.8

# Row 13
$ws.Range("B13").Value = "Experimentación con data aumgentation"
$ws.Range("C13").Value = 13
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 1

# Row 14
$ws.Range("B14").Value = "Experimentacion con oversampling"
$ws.Range("C14").Value = 14
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 14
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 1

# Row 15
$ws.Range("B15").Value = "Comparacion de modelos"
$ws.Range("C15").Value = 15
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 15
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1

# Row 16
$ws.Range("B16").Value = "Entrenamiento con todos los datos"
$ws.Range("C16").Value = 17
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 16
$ws.Range("G16").Value = 0.4

# Row 17
$ws.Range("B17").Value = "Ajustes de arquitectura"
$ws.Range("C17").Value = 18
$ws.Range("D17").Value = 3

# Row 18
$ws.Range("B18").Value = "Pruebas con otras aqrquitecturas"
$ws.Range("C18").Value = 21
$ws.Range("D18").Value = 3

# Row 19
$ws.Range("B19").Value = "Comparacion de resultados"
$ws.Range("C19").Value = 24
$ws.Range("D19").Value = 3

# --- View state: zoom out and move the selection ---
$excel.ActiveWindow.Zoom = 55
$ws.Range("AF24").Select()
